$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the scenario parameter values
$ws.Range("M2").Value = 2
$ws.Range("O2").Value = 2
$ws.Range("R2").Value = 2

# Move the active selection
$ws.Range("C15").Select()
